$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: SCRIPT/T01P01A/um0804.ssb (line 208) and its continuation (line 211).
# Cells are written column-by-column (English, then filename, then Russian
# translation, then converted/cipher text) so the shared-string table grows
# in the same order the source data was produced.
$ws.Range("C3").Value = ' Oh, [hero] and\n[partner]!'
$ws.Range("C4").Value = ' I hope you get lots of treasure!'

$ws.Range("A3").Value = "SCRIPT/T01P01A/um0804.ssb"

$ws.Range("D3").Value = ' Ой, [hero] и\n[partner]!'
$ws.Range("D4").Value = ' Надеюсь, вы найдёте много\nсокровищ!'

$ws.Range("E3").Value = ' Ïê, [hero] é\n[partner]!'
$ws.Range("E4").Value = ' Îàäåýòû, âú îàêäæóå íîïãï\nòïëñïâéþ!'

$ws.Range("B3").Value = 208
$ws.Range("B4").Value = 211

# Row 3 (the row carrying the filename cell) gets the same "two visual
# lines" height used by the other filename rows; row 4 keeps the sheet's
# default height.
$ws.Rows.Item(3).RowHeight = 43.2

# Match the workbook's recorded selection after the edit.
[void]$ws.Range("C2").Select()
